$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 826.86365
$ws.Range("J17").Value = 826.86365
$ws.Range("L17").Value = 2480.59095
$ws.Range("N17").Value = -2816.59095

$ws.Range("H98").Value = 2159.739
$ws.Range("I98").Value = 2246.3809
$ws.Range("J98").Value = 1250
$ws.Range("K98").Value = 2246.3809
$ws.Range("L98").Value = 1250
$ws.Range("M98").Value = -748.3809000000001
$ws.Range("N98").Value = -4246

$ws.Range("H103").Value = 1850
$ws.Range("I103").Value = 1000
$ws.Range("J103").Value = 2133.3333
$ws.Range("K103").Value = 3000
$ws.Range("L103").Value = 6399.999899999999
$ws.Range("M103").Value = -2414
$ws.Range("N103").Value = -7571.999899999999

$ws.Range("H116").Value = 12238.385
$ws.Range("J116").Value = 5149.9
$ws.Range("L116").Value = 5149.9
$ws.Range("N116").Value = -12033.9

$ws.Range("H121").Value = 829.8
$ws.Range("J121").Value = 999.6667
$ws.Range("L121").Value = 2999.0001
$ws.Range("N121").Value = -6493.0001

$ws.Range("H122").Value = 2159.739
$ws.Range("I122").Value = 2246.3809
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 6739.1427
$ws.Range("L122").Value = 3750
$ws.Range("M122").Value = -4289.1427
$ws.Range("N122").Value = -8650

$ws.Range("H137").Value = 30793.734
$ws.Range("I137").Value = 1259.36
$ws.Range("J137").Value = 112833.664
$ws.Range("K137").Value = 3778.08
$ws.Range("L137").Value = 338500.992
$ws.Range("M137").Value = -1228.08
$ws.Range("N137").Value = -343600.992

$ws.Range("H141").Value = 1079619.1
$ws.Range("I141").Value = 1401800.5
$ws.Range("K141").Value = 4205401.5
$ws.Range("M141").Value = -4200221.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 36219.8
$ws.Range("J92").Value = 36219.8
$ws.Range("L92").Value = 36219.8
$ws.Range("N92").Value = -41211.8

$ws.Range("H132").Value = 2385.3333
$ws.Range("J132").Value = 2939.8
$ws.Range("L132").Value = 8819.400000000001
$ws.Range("N132").Value = -13879.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H86").Value = 113654.555
$ws.Range("I86").Value = 2415.4167
$ws.Range("K86").Value = 2415.4167
$ws.Range("M86").Value = -1292.4167

$ws.Range("H89").Value = 113654.555
$ws.Range("I89").Value = 2415.4167
$ws.Range("K89").Value = 12077.0835
$ws.Range("M89").Value = -6461.083500000001

$ws.Range("H105").Value = 2261.652
$ws.Range("J105").Value = 2374.75
$ws.Range("L105").Value = 2374.75
$ws.Range("N105").Value = -5868.75

$ws.Range("H132").Value = 40000
$ws.Range("J132").Value = 40000
$ws.Range("L132").Value = 40000
$ws.Range("N132").Value = -50120

$ws.Range("H134").Value = 7369.875
$ws.Range("I134").Value = 8444.35
$ws.Range("K134").Value = 25333.05
$ws.Range("M134").Value = -22798.05

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 50683.668
$ws.Range("J135").Value = 50683.668
$ws.Range("L135").Value = 50683.668
$ws.Range("N135").Value = -60823.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 295.4
$ws.Range("J23").Value = 295.4
$ws.Range("L23").Value = 886.1999999999999
$ws.Range("N23").Value = -1356.2

$ws.Range("H52").Value = 997.5
$ws.Range("J52").Value = 997.5
$ws.Range("L52").Value = 2992.5
$ws.Range("N52").Value = -3524.5

$ws.Range("H87").Value = 11346.2
$ws.Range("I87").Value = 5851.7144
$ws.Range("K87").Value = 17555.1432
$ws.Range("M87").Value = -16307.1432

$ws.Range("H90").Value = 11346.2
$ws.Range("I90").Value = 5851.7144
$ws.Range("K90").Value = 52665.4296
$ws.Range("M90").Value = -46425.4296

$ws.Range("H105").Value = 2777.7144
$ws.Range("J105").Value = 2923.077
$ws.Range("L105").Value = 8769.231
$ws.Range("N105").Value = -14011.231

$ws.Range("H108").Value = 2002.75
$ws.Range("I108").Value = 2002.75
$ws.Range("K108").Value = 6008.25
$ws.Range("M108").Value = -3128.25

$ws.Range("H109").Value = 1879.5
$ws.Range("I109").Value = 1119.3334
$ws.Range("J109").Value = 4160
$ws.Range("K109").Value = 3358.0002
$ws.Range("L109").Value = 12480
$ws.Range("M109").Value = -2318.0002
$ws.Range("N109").Value = -14560

$ws.Range("H131").Value = 790.25
$ws.Range("J131").Value = 801.30206
$ws.Range("L131").Value = 2403.90618
$ws.Range("N131").Value = -12483.90618

$ws.Range("H137").Value = 2637.3684
$ws.Range("J137").Value = 2853.75
$ws.Range("L137").Value = 8561.25
$ws.Range("N137").Value = -18761.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1854.0625
$ws.Range("I122").Value = 1424.2727
$ws.Range("J122").Value = 2799.6
$ws.Range("K122").Value = 4272.8181
$ws.Range("L122").Value = 8398.799999999999
$ws.Range("M122").Value = -1822.8181
$ws.Range("N122").Value = -13298.8

$ws.Range("H127").Value = 25309.5
$ws.Range("J127").Value = 25309.5
$ws.Range("L127").Value = 25309.5
$ws.Range("N127").Value = -35229.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4186.6665
$ws.Range("I40").Value = 1483.3334
$ws.Range("K40").Value = 1483.3334
$ws.Range("M40").Value = -1347.3334

$ws.Range("H46").Value = 2335.9
$ws.Range("I46").Value = 1415.4
$ws.Range("K46").Value = 1415.4
$ws.Range("M46").Value = -1227.4

$ws.Range("H82").Value = 1954.5883
$ws.Range("I82").Value = 1372.6666
$ws.Range("K82").Value = 1372.6666
$ws.Range("M82").Value = -1011.6666

$ws.Range("H85").Value = 1954.5883
$ws.Range("I85").Value = 1372.6666
$ws.Range("K85").Value = 1372.6666
$ws.Range("M85").Value = -124.6666

$ws.Range("H106").Value = 19793.666
$ws.Range("J106").Value = 19793.666
$ws.Range("L106").Value = 19793.666
$ws.Range("N106").Value = -22317.666

$ws.Range("H122").Value = 12500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 20578820
$ws.Range("J136").Value = 2307.3333
$ws.Range("L136").Value = 6921.999899999999
$ws.Range("N136").Value = -12021.9999
